$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44489
$ws.Range("N2").Value = 26000
$ws.Range("O2").Value = 27000
$ws.Range("P2").Value = 26500
$ws.Range("S2").Value = 1472

# Row 3
$ws.Range("D3").Value = 45043
$ws.Range("L3").Value = 'Segunda'
$ws.Range("M3").Value = 300
$ws.Range("N3").Value = 21000
$ws.Range("O3").Value = 22000
$ws.Range("P3").Value = 21500
$ws.Range("S3").Value = 1194

# Row 4
$ws.Range("D4").Value = 44629
$ws.Range("L4").Value = 'Segunda'
$ws.Range("M4").Value = 300
$ws.Range("N4").Value = 17000
$ws.Range("O4").Value = 18000
$ws.Range("P4").Value = 17500
$ws.Range("S4").Value = 972

# Row 5
$ws.Range("D5").Value = 44819
$ws.Range("M5").Value = 300
$ws.Range("N5").Value = 17000
$ws.Range("O5").Value = 18000
$ws.Range("P5").Value = 17500
$ws.Range("Q5").Value = '$/bandeja 10 kilos'
$ws.Range("S5").Value = 1750
$ws.Range("T5").Value = 10

# Row 6
$ws.Range("D6").Value = 44991
$ws.Range("L6").Value = 'Primera'
$ws.Range("N6").Value = 24000
$ws.Range("O6").Value = 25000
$ws.Range("P6").Value = 24500
$ws.Range("S6").Value = 1361

# Row 8
$ws.Range("D8").Value = 44418
$ws.Range("K8").Value = 'Hayward'
$ws.Range("M8").Value = 240
$ws.Range("N8").Value = 10000
$ws.Range("O8").Value = 11000
$ws.Range("P8").Value = 10500
$ws.Range("Q8").Value = '$/bandeja 10 kilos'
$ws.Range("R8").Value = 'Región de O''Higgins'
$ws.Range("S8").Value = 1050
$ws.Range("T8").Value = 10

# Row 9
$ws.Range("D9").Value = 45034
$ws.Range("M9").Value = 250
$ws.Range("N9").Value = 25000
$ws.Range("O9").Value = 26000
$ws.Range("P9").Value = 25600
$ws.Range("Q9").Value = '$/bandeja 18 kilos'
$ws.Range("S9").Value = 1422
$ws.Range("T9").Value = 18

# Row 10
$ws.Range("D10").Value = 44307
$ws.Range("M10").Value = 250
$ws.Range("N10").Value = 19000
$ws.Range("O10").Value = 20000
$ws.Range("P10").Value = 19500
$ws.Range("Q10").Value = '$/bandeja 18 kilos'
$ws.Range("S10").Value = 1083
$ws.Range("T10").Value = 18

# Row 11
$ws.Range("D11").Value = 44602
$ws.Range("M11").Value = 270
$ws.Range("N11").Value = 20000
$ws.Range("O11").Value = 21000
$ws.Range("P11").Value = 20500
$ws.Range("S11").Value = 1139

# Row 12
$ws.Range("D12").Value = 44263
$ws.Range("M12").Value = 250
$ws.Range("N12").Value = 21000
$ws.Range("O12").Value = 22000
$ws.Range("P12").Value = 21500
$ws.Range("Q12").Value = '$/caja 18 kilos'
$ws.Range("S12").Value = 1194
$ws.Range("T12").Value = 18

# Row 13
$ws.Range("D13").Value = 44784
$ws.Range("M13").Value = 300

# Row 14
$ws.Range("D14").Value = 45002
$ws.Range("L14").Value = 'Segunda'
$ws.Range("M14").Value = 300
$ws.Range("N14").Value = 24000
$ws.Range("O14").Value = 25000
$ws.Range("P14").Value = 24500
$ws.Range("Q14").Value = '$/bandeja 18 kilos'
$ws.Range("S14").Value = 1361
$ws.Range("T14").Value = 18

# Row 15
$ws.Range("D15").Value = 45069
$ws.Range("K15").Value = 'Sin especificar'
$ws.Range("L15").Value = 'Primera'
$ws.Range("M15").Value = 370
$ws.Range("N15").Value = 19000
$ws.Range("O15").Value = 20000
$ws.Range("P15").Value = 19486
$ws.Range("Q15").Value = '$/bandeja 18 kilos'
$ws.Range("R15").Value = 'Región Metropolitana'
$ws.Range("S15").Value = 1083

# Row 16
$ws.Range("D16").Value = 44616
$ws.Range("L16").Value = 'Segunda'
$ws.Range("M16").Value = 300
$ws.Range("N16").Value = 16000
$ws.Range("O16").Value = 17000
$ws.Range("P16").Value = 16500
$ws.Range("Q16").Value = '$/caja 18 kilos granel'
$ws.Range("S16").Value = 917

# Row 17
$ws.Range("D17").Value = 44614
$ws.Range("M17").Value = 250
$ws.Range("N17").Value = 20000
$ws.Range("O17").Value = 21000
$ws.Range("P17").Value = 20500
$ws.Range("Q17").Value = '$/bandeja 18 kilos'
$ws.Range("S17").Value = 1139
$ws.Range("T17").Value = 18

# Row 18
$ws.Range("D18").Value = 44789
$ws.Range("L18").Value = 'Segunda'
$ws.Range("M18").Value = 250

# Row 19
$ws.Range("D19").Value = 44291
$ws.Range("M19").Value = 200
$ws.Range("N19").Value = 17000
$ws.Range("O19").Value = 18000
$ws.Range("P19").Value = 17500
$ws.Range("S19").Value = 972

# Row 20
$ws.Range("D20").Value = 44673
$ws.Range("L20").Value = 'Especial'
$ws.Range("M20").Value = 400
$ws.Range("N20").Value = 14000
$ws.Range("O20").Value = 15000
$ws.Range("P20").Value = 14500
$ws.Range("Q20").Value = '$/bandeja 10 kilos'
$ws.Range("S20").Value = 1450
$ws.Range("T20").Value = 10

# Row 21
$ws.Range("D21").Value = 44487
$ws.Range("M21").Value = 300
$ws.Range("N21").Value = 14000
$ws.Range("O21").Value = 15000
$ws.Range("P21").Value = 14500
$ws.Range("S21").Value = 1450

# Row 22
$ws.Range("D22").Value = 44491
$ws.Range("L22").Value = 'Primera'
$ws.Range("N22").Value = 14000
$ws.Range("O22").Value = 15000
$ws.Range("P22").Value = 14500
$ws.Range("Q22").Value = '$/bandeja 10 kilos'
$ws.Range("S22").Value = 1450
$ws.Range("T22").Value = 10

# Row 23
$ws.Range("D23").Value = 44656
$ws.Range("L23").Value = 'Primera'
$ws.Range("M23").Value = 270
$ws.Range("N23").Value = 19000
$ws.Range("O23").Value = 20000
$ws.Range("P23").Value = 19500
$ws.Range("S23").Value = 1083

# Row 24
$ws.Range("D24").Value = 44706
$ws.Range("M24").Value = 400
$ws.Range("N24").Value = 9000
$ws.Range("O24").Value = 10000
$ws.Range("P24").Value = 9500
$ws.Range("Q24").Value = '$/bandeja 10 kilos'
$ws.Range("S24").Value = 950
$ws.Range("T24").Value = 10
